$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume) hold text-formatted numeric-looking strings.
# Force text format before writing so COM does not coerce them into numbers/dates,
# then restore the default (unstyled) cell style so no stray style index is introduced.
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range('D2').Value = '46.615.32'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.272.64'
$ws.Range('E3').Value = '  -2.55%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '300.81'
$ws.Range('E5').Value = '  -1.96%  '
$ws.Range('D6').Value = '100.52'
$ws.Range('E6').Value = '  +1.81%  '
$ws.Range('D7').Value = '0.568'
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('E8').Value = '  +0.27%  '
$ws.Range('D9').Value = '0.508'
$ws.Range('E9').Value = '  -4.92%  '
$ws.Range('D10').Value = '35.13'
$ws.Range('E10').Value = '  -3.12%  '
$ws.Range('E11').Value = '  -1.25%  '
$ws.Range('D12').Value = '7.08'
$ws.Range('E12').Value = '  -5.15%  '
$ws.Range('E13').Value = '  -1.94%  '
$ws.Range('D14').Value = '2.618.64'
$ws.Range('E14').Value = '  -2.23%  '
$ws.Range('D15').Value = '2.271.18'
$ws.Range('E15').Value = '  -2.15%  '
$ws.Range('E16').Value = '  -3.30%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').Value = '0.800'
$ws.Range('E17').Value = '  -4.17%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '46.601.65'
$ws.Range('E18').Value = '  +0.92%  '
$ws.Range('D19').Value = '12.63'
$ws.Range('E19').Value = '  -3.80%  '
$ws.Range('E20').Value = '  +2.19%  '
$ws.Range('D21').Value = '5.86'
$ws.Range('E21').Value = '  -5.37%  '
$ws.Range('D22').Value = '65.96'
$ws.Range('E22').Value = '  -1.17%  '
$ws.Range('D23').Value = '248.63'
$ws.Range('E23').Value = '  +1.68%  '
$ws.Range('E24').Value = '  -5.40%  '
$ws.Range('E25').Value = '  +0.17%  '
$ws.Range('D26').Value = '1.87'
$ws.Range('E26').Value = '  -5.83%  '
$ws.Range('D27').Value = '41.53'
$ws.Range('E27').Value = '  -1.39%  '
$ws.Range('D29').Value = '9.63'
$ws.Range('E29').Value = '  -2.13%  '
$ws.Range('D30').Value = '20.26'
$ws.Range('E30').Value = '  +0.52%  '
$ws.Range('E31').Value = '  +7.01%  '
$ws.Range('D32').Value = '3.39'
$ws.Range('E32').Value = '  +11.99%  '
$ws.Range('D33').Value = '147.13'
$ws.Range('E33').Value = '  -3.56%  '
$ws.Range('D34').Value = '5.39'
$ws.Range('E34').Value = '  -5.92%  '
$ws.Range('D35').Value = '0.0773'
$ws.Range('E35').Value = '  -4.49%  '
$ws.Range('E36').Value = '  +6.45%  '
$ws.Range('E37').Value = '  -2.72%  '
$ws.Range('D38').Value = '15.81'
$ws.Range('E38').Value = '  +11.09%  '
$ws.Range('D39').Value = '1.69'
$ws.Range('E39').Value = '  -7.48%  '
$ws.Range('D40').Value = '3.88'
$ws.Range('E40').Value = '  -4.48%  '
$ws.Range('E41').Value = '  -6.68%  '
$ws.Range('D42').Value = '3.13'
$ws.Range('E42').Value = '  -7.68%  '
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = '92.54'
$ws.Range('E44').Value = '  +14.74%  '
$ws.Range('D45').Value = '1.790.50'
$ws.Range('E45').Value = '  -0.45%  '
$ws.Range('E46').Value = '  -6.35%  '
$ws.Range('D47').Value = '71.29'
$ws.Range('E47').Value = '  -4.26%  '
$ws.Range('D48').Value = '0.186'
$ws.Range('E48').Value = '  -6.69%  '
$ws.Range('D49').Value = '4.80'
$ws.Range('E49').Value = '  -1.66%  '
$ws.Range('D50').Value = '94.84'
$ws.Range('E50').Value = '  -3.41%  '
$ws.Range('D51').Value = '7.89'
$ws.Range('E51').Value = '  -1.32%  '

$priceVolRange.Style = "Normal"
